$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# --- Sheet1: add a new date column CB ("26-sep") with the day's data ---
$ws1.Range("CB1").Value = "26-sep"

$ws1.Range("CB2").Value  = 0
$ws1.Range("CB3").Value  = 8.6279094617879508
$ws1.Range("CB4").Value  = 10.91905820099767
$ws1.Range("CB5").Value  = 14.712125697480428
$ws1.Range("CB6").Value  = 0
$ws1.Range("CB7").Value  = 0.14341357928414475
$ws1.Range("CB8").Value  = 13.795974844098746
$ws1.Range("CB9").Value  = 10.613387902917527
$ws1.Range("CB10").Value = 13.905110992579916
$ws1.Range("CB11").Value = 8.6126543971001546
$ws1.Range("CB12").Value = 0
$ws1.Range("CB13").Value = 4.0587076097117203
$ws1.Range("CB14").Value = 0
$ws1.Range("CB15").Value = 0
$ws1.Range("CB16").Value = 6.0319720914122206
$ws1.Range("CB17").Value = 0
$ws1.Range("CB18").Value = 0

# Move the active selection like the source session left it.
[void]$ws1.Range("CD10").Select()

# --- Sheet3: correct the typo in the product description (case fix) ---
$ws3.Range("A24").Value = "DORITOS QUESO 85GRX26"

# --- Sheet3: refresh the averages lookup table (A20:B36) with the new day folded in ---
$ws3.Range("B20").Value = 15.230156207210904
$ws3.Range("B21").Value = 4.0587076097117203
$ws3.Range("B22").Value = 13.795974844098746
$ws3.Range("B23").Value = 0.78643886484359637
$ws3.Range("B24").Value = 0.14341357928414475
$ws3.Range("B25").Value = 6.0319720914122206
$ws3.Range("B26").Value = 10.91905820099767
$ws3.Range("B27").Value = 14.712125697480428
$ws3.Range("B28").Value = 3.8699347189568201
$ws3.Range("B29").Value = 8.6279094617879508
$ws3.Range("B30").Value = 16.075829751960246
$ws3.Range("B31").Value = 9.7377703763080934
$ws3.Range("B32").Value = 7.5246468097774271
$ws3.Range("B33").Value = 13.905110992579916
$ws3.Range("B34").Value = 8.6126543971001546
$ws3.Range("B35").Value = 10.613387902917527
$ws3.Range("B36").Value = 26.995786613800714
